$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("H21").Value = "1068 (121)"
$ws.Range("H22").Value = "465.5 (132)"
$ws.Range("H23").Value = "463 (111)"
$ws.Range("H24").Value = "591 (132)"
$ws.Range("H25").Value = "671 (137)"
$ws.Range("H26").Value = "103 (110)"
$ws.Range("H27").Value = "1441.5 (126)"
$ws.Range("H30").Value = "1386 (168)"
$ws.Range("H31").Value = "1575 (112)"
$ws.Range("H32").Value = "590 (125)"
$ws.Range("H33").Value = "301 (132)"
$ws.Range("H34").Value = "1090 (149)"
$ws.Range("H35").Value = "6704 (150)"
$ws.Range("H36").Value = "2390 (131)"
$ws.Range("H38").Value = "1034 (123)"
$ws.Range("H39").Value = "2903 (144)"
$ws.Range("H40").Value = "13703 (114)"
$ws.Range("H43").Value = "1426.5 (128)"
$ws.Range("H44").Value = "12031 (121)"
$ws.Range("H45").Value = "488 (127)"
$ws.Range("H47").Value = "1041 (145)"
$ws.Range("H48").Value = "203 (134)"
$ws.Range("H49").Value = "205 (124)"
$ws.Range("H50").Value = "2032 (138)"
$ws.Range("H51").Value = "903 (125)"
$ws.Range("H53").Value = "1093 (130)"
$ws.Range("H54").Value = "899 (134)"
$ws.Range("H55").Value = "400 (119)"
$ws.Range("H56").Value = "1090 (147)"
$ws.Range("H58").Value = "1993 (144)"
$ws.Range("H59").Value = "600 (141)"
$ws.Range("H60").Value = "280 (164)"
$ws.Range("H61").Value = "1009 (118)"
$ws.Range("H63").Value = "2340 (147)"
$ws.Range("H64").Value = "11930 (125)"
$ws.Range("H65").Value = "14032 (115)"
$ws.Range("H66").Value = "8009 (128)"
$ws.Range("H68").Value = "830 (125)"
$ws.Range("H69").Value = "12930 (120)"
$ws.Range("H70").Value = "5029 (119)"
$ws.Range("H75").Value = "480 (124)"
$ws.Range("H76").Value = "307 (156)"
$ws.Range("H77").Value = "8090 (126)"
$ws.Range("H80").Value = "7361 (113)"
$ws.Range("H82").Value = "339 (119)"
$ws.Range("H83").Value = "1013 (164)"
$ws.Range("H84").Value = "6160 (106)"
$ws.Range("H85").Value = "787 (111)"
$ws.Range("H87").Value = "3322 (138)"
$ws.Range("H89").Value = "712 (111)"
$ws.Range("H90").Value = "392 (108)"
$ws.Range("H92").Value = "1020 (111)"
$ws.Range("H93").Value = "1449 (140)"
$ws.Range("H95").Value = "802 (122)"
$ws.Range("H97").Value = "2009 (100)"
$ws.Range("H99").Value = "3021 (97)"
$ws.Range("H100").Value = "2093 (131)"
$ws.Range("H102").Value = "990 (124)"
$ws.Range("H103").Value = "7092 (86)"
$ws.Range("H104").Value = "1203 (127)"

$ws.Range("R10").Select() | Out-Null
